# revisi buku tesis plus program
#
# Corrects several miscounted "positif"/"negatif" tallies in the monthly /
# per-3-month / per-4-month recap blocks. Every SUM()/percentage formula in
# rows 26, 31, 36, 39, 40, 43, 44, 45 (and the small Q/R summary block)
# depends on these source cells, so simply writing the new literals and
# letting the workbook recalc reproduces the whole cascade of updated
# totals/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24 ("Per bulan" positif row) ---
$ws.Range("B24").Value = 4
$ws.Range("H24").Value = 18
$ws.Range("I24").Value = 10

# --- Row 29 ("Per Tiga Bulan" positif row) ---
$ws.Range("B29").Value = 8
$ws.Range("C29").Value = 7
$ws.Range("D29").Value = 3
$ws.Range("G29").Value = 8
$ws.Range("J29").Value = 10
$ws.Range("K29").Value = 9
$ws.Range("L29").Value = 4
$ws.Range("M29").Value = 1

# --- Row 30 ("Per Tiga Bulan" negatif row) ---
$ws.Range("M30").Value = 4

# --- Row 34 ("Per Empat Bulan" positif row) ---
$ws.Range("M34").Value = 10

# --- Row 35 ("Per Empat Bulan" negatif row) ---
$ws.Range("M35").Value = 7

# --- Leave the sheet scrolled/selected where the author was working ---
try { $excel.ActiveWindow.ScrollRow = 22 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$ws.Range("I25").Select()
